$wb = $excel.ActiveWorkbook

# --- ContactAllocation6_1 (sheet10.xml): insert a new "Department" column ---
$ws1 = $wb.Worksheets.Item("ContactAllocation6_1")
$ws1.Columns.Item(4).Insert() | Out-Null
$ws1.Range("D1").Value = "Department"
$ws1.Range("D2").Value = "QA"
$ws1.Columns.Item(4).ColumnWidth = 10.83

# --- ContactAllocation6_2 (sheet11.xml): insert a new "Department" column ---
$ws2 = $wb.Worksheets.Item("ContactAllocation6_2")
$ws2.Columns.Item(4).Insert() | Out-Null
$ws2.Range("D1").Value = "Department"
$ws2.Range("D2").Value = "QA"
$ws2.Range("D3").Value = "Testing"
$ws2.Range("D4").Value = "Testing"
$ws2.Columns.Item(4).ColumnWidth = 10.83

# --- Selections / active sheet & tab ---
$ws2.Range("A2").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("D2").Select() | Out-Null

Write-Output "done"
